$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
